$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 578, shifting existing rows 578:604 down to 579:605
$ws.Rows("578").Insert()

# Populate the newly inserted row 578 with the new data record
$ws.Range("A578").Value = 3
$ws.Range("B578").Value = "Femacal de La Calera"
$ws.Range("C578").Value = "Coquimbo"
$ws.Range("D578").Value = 45147
$ws.Range("E578").Value = 5
$ws.Range("F578").Value = 100114013
$ws.Range("G578").Value = "Zanahoria"
$ws.Range("H578").Value = "Sin especificar"
$ws.Range("I578").Value = "Primera"
$ws.Range("J578").Value = 260
$ws.Range("K578").Value = 7000
$ws.Range("L578").Value = 7500
$ws.Range("M578").Value = 7308
$ws.Range("N578").Value = "$/saco 20 kilos"
$ws.Range("O578").Value = "Provincia de Quillota"
$ws.Range("P578").Value = 365
$ws.Range("Q578").Value = 20
$ws.Range("R578").Value = "Hortaliza"
